$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 40

# Column A must be stored as text "01/03/2026" (matching the existing date
# cells in the sheet), not auto-converted into a date serial number.
# Temporarily force a text number format so Excel doesn't reinterpret the
# string as a date, then reset the style back to Normal/default so the
# cell doesn't end up with an explicit style index that wasn't in the source.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "01/03/2026"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 12284.42
$ws.Cells.Item($row, 3).Value = 0.2142576612486495
$ws.Cells.Item($row, 4).Value = 0.7857423387513505
$ws.Cells.Item($row, 5).Value = -139.95
$ws.Cells.Item($row, 6).Value = -23.92
$ws.Cells.Item($row, 7).Value = -20989.37
$ws.Cells.Item($row, 8).Value = -68.5
$ws.Cells.Item($row, 9).Value = -420.82
$ws.Cells.Item($row, 10).Value = -13.78
